# Products_Mar2018.xlsx - "Added details related to cap america t shirt"
#
# Sheet3 gets a new product row (Captain America T Shirt) plus a new
# "Status" column with per-row progress, and the selection on Sheet2 /
# Sheet3 is nudged to reflect where the author finished editing.

$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# --- Sheet3: new row for "Captain America T Shirt" -----------------------
# (written first so the shared-string table picks up this text before the
# new "Status" column strings, matching the order the author typed them in)
$ws3.Range("A3").Value = 2
$ws3.Range("B3").Value = "Captain America T Shirt"

# Prices for the existing + new product rows
$ws3.Range("C2").Value = 78.5
$ws3.Range("C3").Value = 81.5

# --- Sheet3: new "Status" column ------------------------------------------
$ws3.Range("D1").Value = "Status"
$ws3.Range("D2").Value = "In Progress"
$ws3.Range("D3").Value = "Completed"

# New (empty-ish) trailing row, only S.No filled in
$ws3.Range("A4").Value = 3

# Column D width to match the other data columns
$ws3.Columns.Item(4).ColumnWidth = 20.6

# --- Selections -------------------------------------------------------
[void]$ws2.Range("B3").Select()
[void]$ws3.Range("A5").Select()
